$p = $ppt.ActivePresentation

# --- Slide 2: add Overview body content, resize content placeholder, ---
# --- expand title text and let it autofit.                            ---
$s2 = $p.Slides.Item(2)

# Content Placeholder 1 (Shapes.Item(1))
$body = $s2.Shapes.Item(1)

# Explicit position/size for the placeholder (EMU -> points, 1 pt = 12700 EMU)
$body.Left   = 457200  / 12700
$body.Top    = 1481328 / 12700
$body.Width  = 8229600 / 12700
$body.Height = 4690872 / 12700

$bodyTextRange = $body.TextFrame.TextRange
$lines = @(
    "Structure: Face page, approval, contents",
    "Description in natural language",
    "3 classes: Activator, Driver, Message Handler",
    "Activator makes connection, sends and receives messages, and has instances of Driver and Message Handler",
    "Driver controls robot movement",
    "Message Handler decodes, encodes, and verifies messages",
    "UML Class Diagram",
    "For postlab: UML sequence diagram",
    ""
)
$bodyTextRange.Text = [string]::Join("`r", $lines)

# Outline / indent levels (COM IndentLevel is 1-based; OOXML lvl = IndentLevel-1)
$bodyTextRange.Paragraphs(1, 1).IndentLevel = 1
$bodyTextRange.Paragraphs(2, 1).IndentLevel = 1
$bodyTextRange.Paragraphs(3, 1).IndentLevel = 2
$bodyTextRange.Paragraphs(4, 1).IndentLevel = 2
$bodyTextRange.Paragraphs(5, 1).IndentLevel = 2
$bodyTextRange.Paragraphs(6, 1).IndentLevel = 2
$bodyTextRange.Paragraphs(7, 1).IndentLevel = 1
$bodyTextRange.Paragraphs(8, 1).IndentLevel = 1
$bodyTextRange.Paragraphs(9, 1).IndentLevel = 3

# Title 2 (Shapes.Item(2)) -- longer title, shrink-to-fit turned on
$title = $s2.Shapes.Item(2)
$title.TextFrame.TextRange.Text = "Overview of Design Documentation"
$title.TextFrame.AutoSize = 2

Write-Host "Slide 2 updated"
